# Add an "Average" summary row (row 5) to every worksheet in the workbook.
# Row 5, column A gets the label "Average"; every data column that has
# values in rows 2-4 gets "=AVERAGE(<col>2:<col>4)" in row 5.

$wb = $excel.ActiveWorkbook

$cols = @("B","E","H","K","N","Q","T","W","Z","AC","AF","AI","AL","AO","AR")

foreach ($ws in $wb.Worksheets) {
    $ws.Range("A5").Value = "Average"
    foreach ($col in $cols) {
        $cellRef = "${col}5"
        $ws.Range($cellRef).Formula = "=AVERAGE(${col}2:${col}4)"
    }
}
